$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text, since values like "67.339.12" or "580.47"
# are not valid numbers / would otherwise be auto-converted by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.307.02'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").Value = '3.136.11'
$ws.Range("E3").Value = '  +3.48%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '580.47'
$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("D6").Value = '174.90'
$ws.Range("E6").Value = '  +3.85%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.130.79'
$ws.Range("E8").Value = '  +3.43%  '

$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("E10").Value = '  -2.26%  '

$ws.Range("E11").Value = '  +2.16%  '

$ws.Range("D12").Value = '0.484'
$ws.Range("E12").Value = '  -1.13%  '

$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").Value = '37.44'
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("E15").Value = '  -0.76%  '

$ws.Range("D16").Value = '3.657.95'
$ws.Range("E16").Value = '  +3.55%  '

$ws.Range("D17").Value = '67.304.66'
$ws.Range("E17").Value = '  +1.57%  '

$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").Value = '3.140.29'
$ws.Range("E19").Value = '  +3.72%  '

$ws.Range("E20").Value = '  -1.84%  '

$ws.Range("D21").Value = '488.19'
$ws.Range("E21").Value = '  +4.17%  '

$ws.Range("D22").Value = '0.717'
$ws.Range("E22").Value = '  +1.05%  '

$ws.Range("D23").Value = '7.71'
$ws.Range("E23").Value = '  +4.35%  '

$ws.Range("D24").Value = '84.31'
$ws.Range("E24").Value = '  +1.51%  '

$ws.Range("D25").Value = '13.25'
$ws.Range("E25").Value = '  +3.80%  '

$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +3.11%  '

$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  +0.55%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").Value = '7.98'
$ws.Range("E29").Value = '  -2.81%  '

$ws.Range("E30").Value = '  -0.59%  '

$ws.Range("D31").Value = '2.69'
$ws.Range("E31").Value = '  +1.67%  '

$ws.Range("D32").Value = '28.86'
$ws.Range("E32").Value = '  +2.33%  '

$ws.Range("E33").Value = '  +0.31%  '

$ws.Range("E34").Value = '  -3.44%  '

$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").Value = '5.95'
$ws.Range("E36").Value = '  +1.36%  '

$ws.Range("D37").Value = '0.988'
$ws.Range("E37").Value = '  -0.36%  '

$ws.Range("D38").Value = '47.56'
$ws.Range("E38").Value = '  -1.39%  '

$ws.Range("D39").Value = '2.11'
$ws.Range("E39").Value = '  +2.25%  '

$ws.Range("D40").Value = '50.14'
$ws.Range("E40").Value = '  +1.31%  '

$ws.Range("D41").Value = '0.313'
$ws.Range("E41").Value = '  +0.34%  '

$ws.Range("E42").Value = '  +1.71%  '

$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").Value = '2.81'
$ws.Range("E44").Value = '  -1.44%  '

$ws.Range("D45").Value = '2.850.77'
$ws.Range("E45").Value = '  +5.29%  '

$ws.Range("D46").Value = '385.26'
$ws.Range("E46").Value = '  +1.69%  '

$ws.Range("E47").Value = '  -0.50%  '

$ws.Range("D48").Value = '135.98'
$ws.Range("E48").Value = '  +1.22%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").Value = '24.98'
$ws.Range("E50").Value = '  +1.90%  '

$ws.Range("E51").Value = '  -0.38%  '
